$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.424164
$ws.Range("H2").Value = 4.272492
$ws.Range("I2").Value = 0.03823384556371837
$ws.Range("J2").Value = 0.03823384556371837
$ws.Range("M2").Value = 0.1005383333333333
$ws.Range("N2").Value = 0.301615
$ws.Range("O2").Value = 0.0006585930795375908
$ws.Range("P2").Value = 0.0006585930795375907
$ws.Range("Q2").Value = 0.1431830749533334
$ws.Range("R2").Value = 1.28864767458
$ws.Range("S2").Value = 0.000025180546092373940115670058
$ws.Range("T2").Value = 0.000025180546092373929951274691

$ws.Range("G3").Value = 1.424164
$ws.Range("H3").Value = 4.272492
$ws.Range("I3").Value = 0.03823384556371837
$ws.Range("J3").Value = 0.03823384556371837
$ws.Range("M3").Value = 89.19784566666665
$ws.Range("N3").Value = 267.593537
$ws.Range("O3").Value = 0.5843053283065703
$ws.Range("P3").Value = 0.5843053283065702
$ws.Range("Q3").Value = 127.0323606760226
$ws.Range("R3").Value = 1143.291246084204
$ws.Range("S3").Value = 0.02234023968453117
$ws.Range("T3").Value = 0.02234023968453116

$ws.Range("G4").Value = 1.424164
$ws.Range("H4").Value = 4.272492
$ws.Range("I4").Value = 0.03823384556371837
$ws.Range("J4").Value = 0.03823384556371837
$ws.Range("O4").Value = 0.0003271664814192805
$ws.Range("P4").Value = 0.0003271664814192805
$ws.Range("Q4").Value = 0.07112844681599999
$ws.Range("R4").Value = 0.640156021344
$ws.Range("S4").Value = 0.000012508832724209909274988572
$ws.Range("T4").Value = 0.0000125088327242099008046591

$ws.Range("G5").Value = 1.424164
$ws.Range("H5").Value = 4.272492
$ws.Range("I5").Value = 0.03823384556371837
$ws.Range("J5").Value = 0.03823384556371837
$ws.Range("M5").Value = 63.307897
$ws.Range("N5").Value = 189.923691
$ws.Range("O5").Value = 0.4147089121324729
$ws.Range("P5").Value = 0.4147089121324729
$ws.Range("Q5").Value = 90.16082782310801
$ws.Range("R5").Value = 811.447450407972
$ws.Range("S5").Value = 0.01585591650037062
$ws.Range("T5").Value = 0.01585591650037062

$ws.Range("I6").Value = 0.08783524098133262
$ws.Range("J6").Value = 0.08783524098133261
$ws.Range("M6").Value = 0.1005383333333333
$ws.Range("N6").Value = 0.301615
$ws.Range("O6").Value = 0.0006585930795375908
$ws.Range("P6").Value = 0.0006585930795375907
$ws.Range("Q6").Value = 0.3289368282877778
$ws.Range("R6").Value = 2.96043145459
$ws.Range("S6").Value = 0.000057847681849822246886744931
$ws.Range("T6").Value = 0.000057847681849822233334217775

$ws.Range("I7").Value = 0.08783524098133262
$ws.Range("J7").Value = 0.08783524098133261
$ws.Range("M7").Value = 89.19784566666665
$ws.Range("N7").Value = 267.593537
$ws.Range("O7").Value = 0.5843053283065703
$ws.Range("P7").Value = 0.5843053283065702
$ws.Range("Q7").Value = 291.8335272817602
$ws.Range("R7").Value = 2626.501745535842
$ws.Range("S7").Value = 0.05132259931848428
$ws.Range("T7").Value = 0.05132259931848426

$ws.Range("I8").Value = 0.08783524098133262
$ws.Range("J8").Value = 0.08783524098133261
$ws.Range("O8").Value = 0.0003271664814192805
$ws.Range("P8").Value = 0.0003271664814192805
$ws.Range("S8").Value = 0.000028736746736477188502844465
$ws.Range("T8").Value = 0.00002873674673647717156218552

$ws.Range("I9").Value = 0.08783524098133262
$ws.Range("J9").Value = 0.08783524098133261
$ws.Range("M9").Value = 63.307897
$ws.Range("N9").Value = 189.923691
$ws.Range("O9").Value = 0.4147089121324729
$ws.Range("P9").Value = 0.4147089121324729
$ws.Range("Q9").Value = 207.1279496518673
$ws.Range("R9").Value = 1864.151546866806
$ws.Range("S9").Value = 0.03642605723426205
$ws.Range("T9").Value = 0.03642605723426205

$ws.Range("G10").Value = 17.63507366666667
$ws.Range("H10").Value = 52.905221
$ws.Range("I10").Value = 0.4734403362787783
$ws.Range("J10").Value = 0.4734403362787782
$ws.Range("M10").Value = 0.1005383333333333
$ws.Range("N10").Value = 0.301615
$ws.Range("O10").Value = 0.0006585930795375908
$ws.Range("P10").Value = 0.0006585930795375907
$ws.Range("Q10").Value = 1.773000914657222
$ws.Range("R10").Value = 15.957008231915
$ws.Range("S10").Value = 0.0003118045290471531
$ws.Range("T10").Value = 0.0003118045290471531

$ws.Range("G11").Value = 17.63507366666667
$ws.Range("H11").Value = 52.905221
$ws.Range("I11").Value = 0.4734403362787783
$ws.Range("J11").Value = 0.4734403362787782
$ws.Range("M11").Value = 89.19784566666665
$ws.Range("N11").Value = 267.593537
$ws.Range("O11").Value = 0.5843053283065703
$ws.Range("P11").Value = 0.5843053283065702
$ws.Range("Q11").Value = 1573.01057923963
$ws.Range("R11").Value = 14157.09521315668
$ws.Range("S11").Value = 0.2766337111229446
$ws.Range("T11").Value = 0.2766337111229445

$ws.Range("G12").Value = 17.63507366666667
$ws.Range("H12").Value = 52.905221
$ws.Range("I12").Value = 0.4734403362787783
$ws.Range("J12").Value = 0.4734403362787782
$ws.Range("O12").Value = 0.0003271664814192805
$ws.Range("P12").Value = 0.0003271664814192805
$ws.Range("Q12").Value = 0.880766119208
$ws.Range("R12").Value = 7.926895072871999
$ws.Range("S12").Value = 0.0001548938089822888
$ws.Range("T12").Value = 0.0001548938089822888

$ws.Range("G13").Value = 17.63507366666667
$ws.Range("H13").Value = 52.905221
$ws.Range("I13").Value = 0.4734403362787783
$ws.Range("J13").Value = 0.4734403362787782
$ws.Range("M13").Value = 63.307897
$ws.Range("N13").Value = 189.923691
$ws.Range("O13").Value = 0.4147089121324729
$ws.Range("P13").Value = 0.4147089121324729
$ws.Range("Q13").Value = 1116.439427276746
$ws.Range("R13").Value = 10047.95484549071
$ws.Range("S13").Value = 0.1963399268178043
$ws.Range("T13").Value = 0.1963399268178042

$ws.Range("G14").Value = 0.5460243333333333
$ws.Range("H14").Value = 1.638073
$ws.Range("I14").Value = 0.01465885251607185
$ws.Range("J14").Value = 0.01465885251607185
$ws.Range("M14").Value = 0.1005383333333333
$ws.Range("N14").Value = 0.301615
$ws.Range("O14").Value = 0.0006585930795375908
$ws.Range("P14").Value = 0.0006585930795375907
$ws.Range("Q14").Value = 0.05489637643277778
$ws.Range("R14").Value = 0.494067387895
$ws.Range("S14").Value = 0.000009654218821047118979255565
$ws.Range("T14").Value = 0.00000965421882104711728518967

$ws.Range("G15").Value = 0.5460243333333333
$ws.Range("H15").Value = 1.638073
$ws.Range("I15").Value = 0.01465885251607185
$ws.Range("J15").Value = 0.01465885251607185
$ws.Range("M15").Value = 89.19784566666665
$ws.Range("N15").Value = 267.593537
$ws.Range("O15").Value = 0.5843053283065703
$ws.Range("P15").Value = 0.5843053283065702
$ws.Range("Q15").Value = 48.70419421491121
$ws.Range("R15").Value = 438.337747934201
$ws.Range("S15").Value = 0.008565245632000956
$ws.Range("T15").Value = 0.008565245632000954

$ws.Range("G16").Value = 0.5460243333333333
$ws.Range("H16").Value = 1.638073
$ws.Range("I16").Value = 0.01465885251607185
$ws.Range("J16").Value = 0.01465885251607185
$ws.Range("O16").Value = 0.0003271664814192805
$ws.Range("P16").Value = 0.0003271664814192805
$ws.Range("Q16").Value = 0.027270639304
$ws.Range("R16").Value = 0.245435753736
$ws.Range("S16").Value = 0.000004795885199327393792680541
$ws.Range("T16").Value = 0.000004795885199327392945647594

$ws.Range("G17").Value = 0.5460243333333333
$ws.Range("H17").Value = 1.638073
$ws.Range("I17").Value = 0.01465885251607185
$ws.Range("J17").Value = 0.01465885251607185
$ws.Range("M17").Value = 63.307897
$ws.Range("N17").Value = 189.923691
$ws.Range("O17").Value = 0.4147089121324729
$ws.Range("P17").Value = 0.4147089121324729
$ws.Range("Q17").Value = 34.56765225416034
$ws.Range("R17").Value = 311.108870287443
$ws.Range("S17").Value = 0.006079156780050519
$ws.Range("T17").Value = 0.006079156780050519

$ws.Range("G18").Value = 14.37176
$ws.Range("H18").Value = 43.11528
$ws.Range("I18").Value = 0.385831724660099
$ws.Range("J18").Value = 0.3858317246600989
$ws.Range("M18").Value = 0.1005383333333333
$ws.Range("N18").Value = 0.301615
$ws.Range("O18").Value = 0.0006585930795375908
$ws.Range("P18").Value = 0.0006585930795375907
$ws.Range("Q18").Value = 1.444912797466667
$ws.Range("R18").Value = 13.0042151772
$ws.Range("S18").Value = 0.0002541061037271944
$ws.Range("T18").Value = 0.0002541061037271943

$ws.Range("G19").Value = 14.37176
$ws.Range("H19").Value = 43.11528
$ws.Range("I19").Value = 0.385831724660099
$ws.Range("J19").Value = 0.3858317246600989
$ws.Range("M19").Value = 89.19784566666665
$ws.Range("N19").Value = 267.593537
$ws.Range("O19").Value = 0.5843053283065703
$ws.Range("P19").Value = 0.5843053283065702
$ws.Range("Q19").Value = 1281.930030438373
$ws.Range("R19").Value = 11537.37027394536
$ws.Range("S19").Value = 0.2254435325486094
$ws.Range("T19").Value = 0.2254435325486093

$ws.Range("G20").Value = 14.37176
$ws.Range("H20").Value = 43.11528
$ws.Range("I20").Value = 0.385831724660099
$ws.Range("J20").Value = 0.3858317246600989
$ws.Range("O20").Value = 0.0003271664814192805
$ws.Range("P20").Value = 0.0003271664814192805
$ws.Range("Q20").Value = 0.71778318144
$ws.Range("R20").Value = 6.46004863296
$ws.Range("S20").Value = 0.0001262312077769772
$ws.Range("T20").Value = 0.0001262312077769772

$ws.Range("G21").Value = 14.37176
$ws.Range("H21").Value = 43.11528
$ws.Range("I21").Value = 0.385831724660099
$ws.Range("J21").Value = 0.3858317246600989
$ws.Range("M21").Value = 63.307897
$ws.Range("N21").Value = 189.923691
$ws.Range("O21").Value = 0.4147089121324729
$ws.Range("P21").Value = 0.4147089121324729
$ws.Range("Q21").Value = 909.8459017887201
$ws.Range("R21").Value = 8188.613116098481
$ws.Range("S21").Value = 0.1600078547999855
$ws.Range("T21").Value = 0.1600078547999854
